$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 171
$ws.Range("I2").Value = 541
$ws.Range("J2").Value = 2211
$ws.Range("L2").Value = 584
$ws.Range("M2").Value = 33
$ws.Range("N2").Value = 410
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 20
$ws.Range("S2").Value = 244
$ws.Range("T2").Value = 392
$ws.Range("U2").Value = 27
$ws.Range("V2").Value = 3345
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 3431
$ws.Range("Z2").Value = 51
$ws.Range("AA2").Value = 20
